$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '98.471.15'
$ws.Range('E2').Value = '  +0.10%  '
$ws.Range('D3').Value = '3.354.35'
$ws.Range('E3').Value = '  +0.85%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '257.35'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.08%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '663.73'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +6.18%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.53'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +9.15%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.475'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +22.60%  '
$ws.Range('E9').Value = '  +24.43%  '
$ws.Range('E10').Value = '  -0.10%  '
$ws.Range('D11').Value = '3.351.07'
$ws.Range('E11').Value = '  +0.90%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.216'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +8.73%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '42.22'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +14.38%  '
$ws.Range('E14').Value = '  +10.69%  '
$ws.Range('D15').Value = '98.333.29'
$ws.Range('E15').Value = '  +0.14%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.70'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +3.79%  '
$ws.Range('D17').Value = '3.981.40'
$ws.Range('E17').Value = '  +1.08%  '
$ws.Range('D18').Value = '3.355.64'
$ws.Range('E18').Value = '  +0.72%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.61'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +25.91%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '16.69'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +10.78%  '
$ws.Range('B21').Value = 'SuiNetwork'
$ws.Range('C21').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '3.57'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +1.34%  '
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '528.96'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +8.62%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '10.62'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +13.92%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.0000219'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +4.28%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.435'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +52.99%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '102.18'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +15.39%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '6.18'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +10.35%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '12.51'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +5.44%  '
$ws.Range('D29').Value = '3.533.69'
$ws.Range('E29').Value = '  +0.79%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.148'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +8.02%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '11.01'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +14.96%  '
$ws.Range('E33').Value = '  -1.78%  '
$ws.Range('E34').Value = '  +0.28%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '29.28'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +5.71%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.538'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +17.82%  '
$ws.Range('E37').Value = '  +7.21%  '
$ws.Range('E38').Value = '  +8.90%  '
$ws.Range('E39').Value = '  +5.67%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '524.85'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +6.22%  '
$ws.Range('E41').Value = '  -0.58%  '
$ws.Range('E42').Value = '  +6.27%  '
$ws.Range('E43').Value = '  +4.23%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0430'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +32.71%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.42'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +4.03%  '
$ws.Range('E46').Value = '  +5.69%  '
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '5.13'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +11.82%  '
$ws.Range('E49').Value = '  +7.00%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.83'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +17.39%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '50.80'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +11.71%  '
